# Update NATMI LR-pair output with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 29.223446
$ws.Cells.Item(2, 8).Value = 87.670338
$ws.Cells.Item(2, 9).Value = 0.0169041244192178
$ws.Cells.Item(2, 10).Value = 0.0169041244192178
$ws.Cells.Item(2, 13).Value = 3.425446666666666
$ws.Cells.Item(2, 14).Value = 10.27634
$ws.Cells.Item(2, 15).Value = 0.6657953389778073
$ws.Cells.Item(2, 16).Value = 0.6657953389778073
$ws.Cells.Item(2, 17).Value = 100.1033556892133
$ws.Cells.Item(2, 18).Value = 900.93020120292
$ws.Cells.Item(2, 19).Value = 0.01125468724781615
$ws.Cells.Item(2, 20).Value = 0.01125468724781615

$ws.Cells.Item(3, 7).Value = 29.223446
$ws.Cells.Item(3, 8).Value = 87.670338
$ws.Cells.Item(3, 9).Value = 0.0169041244192178
$ws.Cells.Item(3, 10).Value = 0.0169041244192178
$ws.Cells.Item(3, 15).Value = 0.2094791321596951
$ws.Cells.Item(3, 16).Value = 0.2094791321596952
$ws.Cells.Item(3, 17).Value = 31.495510479608
$ws.Cells.Item(3, 18).Value = 283.459594316472
$ws.Cells.Item(3, 19).Value = 0.003541061313257256
$ws.Cells.Item(3, 20).Value = 0.003541061313257256

$ws.Cells.Item(4, 7).Value = 29.223446
$ws.Cells.Item(4, 8).Value = 87.670338
$ws.Cells.Item(4, 9).Value = 0.0169041244192178
$ws.Cells.Item(4, 10).Value = 0.0169041244192178
$ws.Cells.Item(4, 13).Value = 0.62317
$ws.Cells.Item(4, 14).Value = 1.86951
$ws.Cells.Item(4, 15).Value = 0.1211239647746572
$ws.Cells.Item(4, 16).Value = 0.1211239647746572
$ws.Cells.Item(4, 17).Value = 18.21117484382
$ws.Cells.Item(4, 18).Value = 163.90057359438
$ws.Cells.Item(4, 19).Value = 0.00204749457069976
$ws.Cells.Item(4, 20).Value = 0.00204749457069976

$ws.Cells.Item(5, 7).Value = 29.223446
$ws.Cells.Item(5, 8).Value = 87.670338
$ws.Cells.Item(5, 9).Value = 0.0169041244192178
$ws.Cells.Item(5, 10).Value = 0.0169041244192178
$ws.Cells.Item(5, 13).Value = 0.01852966666666667
$ws.Cells.Item(5, 14).Value = 0.055589
$ws.Cells.Item(5, 15).Value = 0.003601564087840353
$ws.Cells.Item(5, 16).Value = 0.003601564087840353
$ws.Cells.Item(5, 17).Value = 0.5415007132313333
$ws.Cells.Item(5, 18).Value = 4.873506419082
$ws.Cells.Item(5, 19).Value = 0.00006088128744464
$ws.Cells.Item(5, 20).Value = 0.00006088128744464

$ws.Cells.Item(6, 9).Value = 0.9471112884046843
$ws.Cells.Item(6, 10).Value = 0.9471112884046842
$ws.Cells.Item(6, 13).Value = 3.425446666666666
$ws.Cells.Item(6, 14).Value = 10.27634
$ws.Cells.Item(6, 15).Value = 0.6657953389778073
$ws.Cells.Item(6, 16).Value = 0.6657953389778073
$ws.Cells.Item(6, 17).Value = 5608.632297610021
$ws.Cells.Item(6, 18).Value = 50477.69067849019
$ws.Cells.Item(6, 19).Value = 0.6305822813131046
$ws.Cells.Item(6, 20).Value = 0.6305822813131045

$ws.Cells.Item(7, 9).Value = 0.9471112884046843
$ws.Cells.Item(7, 10).Value = 0.9471112884046842
$ws.Cells.Item(7, 15).Value = 0.2094791321596951
$ws.Cells.Item(7, 16).Value = 0.2094791321596952
$ws.Cells.Item(7, 19).Value = 0.198400050753664
$ws.Cells.Item(7, 20).Value = 0.198400050753664

$ws.Cells.Item(8, 9).Value = 0.9471112884046843
$ws.Cells.Item(8, 10).Value = 0.9471112884046842
$ws.Cells.Item(8, 13).Value = 0.62317
$ws.Cells.Item(8, 14).Value = 1.86951
$ws.Cells.Item(8, 15).Value = 0.1211239647746572
$ws.Cells.Item(8, 16).Value = 0.1211239647746572
$ws.Cells.Item(8, 17).Value = 1020.343251265033
$ws.Cells.Item(8, 18).Value = 9183.089261385299
$ws.Cells.Item(8, 19).Value = 0.1147178743344092
$ws.Cells.Item(8, 20).Value = 0.1147178743344092

$ws.Cells.Item(9, 9).Value = 0.9471112884046843
$ws.Cells.Item(9, 10).Value = 0.9471112884046842
$ws.Cells.Item(9, 13).Value = 0.01852966666666667
$ws.Cells.Item(9, 14).Value = 0.055589
$ws.Cells.Item(9, 15).Value = 0.003601564087840353
$ws.Cells.Item(9, 16).Value = 0.003601564087840353
$ws.Cells.Item(9, 17).Value = 30.33942637085222
$ws.Cells.Item(9, 18).Value = 273.05483733767
$ws.Cells.Item(9, 19).Value = 0.003411082003506518
$ws.Cells.Item(9, 20).Value = 0.003411082003506518

$ws.Cells.Item(10, 7).Value = 37.39212666666667
$ws.Cells.Item(10, 8).Value = 112.17638
$ws.Cells.Item(10, 9).Value = 0.02162924801792661
$ws.Cells.Item(10, 10).Value = 0.0216292480179266
$ws.Cells.Item(10, 13).Value = 3.425446666666666
$ws.Cells.Item(10, 14).Value = 10.27634
$ws.Cells.Item(10, 15).Value = 0.6657953389778073
$ws.Cells.Item(10, 16).Value = 0.6657953389778073
$ws.Cells.Item(10, 17).Value = 128.0847356499111
$ws.Cells.Item(10, 18).Value = 1152.7626208492
$ws.Cells.Item(10, 19).Value = 0.01440065251593051
$ws.Cells.Item(10, 20).Value = 0.01440065251593051

$ws.Cells.Item(11, 7).Value = 37.39212666666667
$ws.Cells.Item(11, 8).Value = 112.17638
$ws.Cells.Item(11, 9).Value = 0.02162924801792661
$ws.Cells.Item(11, 10).Value = 0.0216292480179266
$ws.Cells.Item(11, 15).Value = 0.2094791321596951
$ws.Cells.Item(11, 16).Value = 0.2094791321596952
$ws.Cells.Item(11, 17).Value = 40.29928973074666
$ws.Cells.Item(11, 18).Value = 362.69360757672
$ws.Cells.Item(11, 19).Value = 0.004530876104062072
$ws.Cells.Item(11, 20).Value = 0.004530876104062072

$ws.Cells.Item(12, 7).Value = 37.39212666666667
$ws.Cells.Item(12, 8).Value = 112.17638
$ws.Cells.Item(12, 9).Value = 0.02162924801792661
$ws.Cells.Item(12, 10).Value = 0.0216292480179266
$ws.Cells.Item(12, 13).Value = 0.62317
$ws.Cells.Item(12, 14).Value = 1.86951
$ws.Cells.Item(12, 15).Value = 0.1211239647746572
$ws.Cells.Item(12, 16).Value = 0.1211239647746572
$ws.Cells.Item(12, 17).Value = 23.30165157486667
$ws.Cells.Item(12, 18).Value = 209.7148641738
$ws.Cells.Item(12, 19).Value = 0.002619820275025666
$ws.Cells.Item(12, 20).Value = 0.002619820275025666

$ws.Cells.Item(13, 7).Value = 37.39212666666667
$ws.Cells.Item(13, 8).Value = 112.17638
$ws.Cells.Item(13, 9).Value = 0.02162924801792661
$ws.Cells.Item(13, 10).Value = 0.0216292480179266
$ws.Cells.Item(13, 13).Value = 0.01852966666666667
$ws.Cells.Item(13, 14).Value = 0.055589
$ws.Cells.Item(13, 15).Value = 0.003601564087840353
$ws.Cells.Item(13, 16).Value = 0.003601564087840353
$ws.Cells.Item(13, 17).Value = 0.6928636430911111
$ws.Cells.Item(13, 18).Value = 6.23577278782
$ws.Cells.Item(13, 19).Value = 0.00007789912290835661
$ws.Cells.Item(13, 20).Value = 0.0000778991229083566

$ws.Cells.Item(14, 7).Value = 24.817167
$ws.Cells.Item(14, 8).Value = 74.45150100000001
$ws.Cells.Item(14, 9).Value = 0.01435533915817136
$ws.Cells.Item(14, 10).Value = 0.01435533915817136
$ws.Cells.Item(14, 13).Value = 3.425446666666666
$ws.Cells.Item(14, 14).Value = 10.27634
$ws.Cells.Item(14, 15).Value = 0.6657953389778073
$ws.Cells.Item(14, 16).Value = 0.6657953389778073
$ws.Cells.Item(14, 17).Value = 85.00988197625999
$ws.Cells.Item(14, 18).Value = 765.08893778634
$ws.Cells.Item(14, 19).Value = 0.009557717900956093
$ws.Cells.Item(14, 20).Value = 0.009557717900956091

$ws.Cells.Item(15, 7).Value = 24.817167
$ws.Cells.Item(15, 8).Value = 74.45150100000001
$ws.Cells.Item(15, 9).Value = 0.01435533915817136
$ws.Cells.Item(15, 10).Value = 0.01435533915817136
$ws.Cells.Item(15, 15).Value = 0.2094791321596951
$ws.Cells.Item(15, 16).Value = 0.2094791321596952
$ws.Cells.Item(15, 17).Value = 26.746652099916
$ws.Cells.Item(15, 18).Value = 240.719868899244
$ws.Cells.Item(15, 19).Value = 0.003007143988711825
$ws.Cells.Item(15, 20).Value = 0.003007143988711825

$ws.Cells.Item(16, 7).Value = 24.817167
$ws.Cells.Item(16, 8).Value = 74.45150100000001
$ws.Cells.Item(16, 9).Value = 0.01435533915817136
$ws.Cells.Item(16, 10).Value = 0.01435533915817136
$ws.Cells.Item(16, 13).Value = 0.62317
$ws.Cells.Item(16, 14).Value = 1.86951
$ws.Cells.Item(16, 15).Value = 0.1211239647746572
$ws.Cells.Item(16, 16).Value = 0.1211239647746572
$ws.Cells.Item(16, 17).Value = 15.46531395939
$ws.Cells.Item(16, 18).Value = 139.18782563451
$ws.Cells.Item(16, 19).Value = 0.001738775594522605
$ws.Cells.Item(16, 20).Value = 0.001738775594522605

$ws.Cells.Item(17, 7).Value = 24.817167
$ws.Cells.Item(17, 8).Value = 74.45150100000001
$ws.Cells.Item(17, 9).Value = 0.01435533915817136
$ws.Cells.Item(17, 10).Value = 0.01435533915817136
$ws.Cells.Item(17, 13).Value = 0.01852966666666667
$ws.Cells.Item(17, 14).Value = 0.055589
$ws.Cells.Item(17, 15).Value = 0.003601564087840353
$ws.Cells.Item(17, 16).Value = 0.003601564087840353
$ws.Cells.Item(17, 17).Value = 0.459853832121
$ws.Cells.Item(17, 18).Value = 4.138684489089
$ws.Cells.Item(17, 19).Value = 0.00005170167398083833
$ws.Cells.Item(17, 20).Value = 0.00005170167398083833
